$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (column C) and P_Value (column D) for rows 2-11
$ws.Range("C2").Value = -0.1183573370942024
$ws.Range("D2").Value = 0.906481031471261

$ws.Range("C3").Value = 0.07613461832351974
$ws.Range("D3").Value = 0.9397582409044707

$ws.Range("C4").Value = 3.626171437064089
$ws.Range("D4").Value = 0.0009316760458732976

$ws.Range("C5").Value = 3.334993276412798
$ws.Range("D5").Value = 0.002071206097957745

$ws.Range("C6").Value = 0.2240742811785195
$ws.Range("D6").Value = 0.8240405718134338

$ws.Range("C7").Value = 5.087106460513424
$ws.Range("D7").Value = 0.00001324391514123668

$ws.Range("C8").Value = 4.103493458588746
$ws.Range("D8").Value = 0.0002404574569832008

$ws.Range("C9").Value = 4.785891634404568
$ws.Range("D9").Value = 0.00003249406826966172

$ws.Range("C10").Value = 5.878888308529246
$ws.Range("D10").Value = 0.000001236973069751457

$ws.Range("C11").Value = -1.249362027993026
$ws.Range("D11").Value = 0.2200655082014473
